$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the RS-485 transceiver part back to THVD2410 (was SN65HVD72DR).
# Set the text, then re-apply the original cell formatting (PasteSpecial
# formats only) so the style stays consistent with the rest of column A.
$ws.Range("A8").Value = "THVD2410"
$ws.Range("B8").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the SM712.TCT (U2, TVS diode) row entirely - no longer used.
$ws.Rows(9).Delete()
